$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the "closing" (bottom-border) formatting from the last table row (25)
#        onto row 18, which becomes the new last row of the (now 3-row) table once
#        the rows below it are removed.
$ws.Range("B25:J25").Copy()
$ws.Range("B18:J18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# --- 2. Re-order the "Periodo Mora" values for LUIS PUENTES VALLE from descending
#        (1803,1802,1801) to ascending (1801,1802,1803).
$ws.Range("E16").Value = "1801"
$ws.Range("E17").Value = "1802"
$ws.Range("E18").Value = "1803"

# --- 3. Update the summary header values.
$ws.Range("E11").Value = 120000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 3

# --- 4. Remove the CARLOS ANDRES ALVAREZ TAMARA block of rows (old rows 19-25),
#        which also shifts the footer/signature rows up into place.
$ws.Rows("19:25").Delete()
